# Convert OOXML EMU values to points (PowerPoint COM uses points for Shape coordinates).
# 1 point = 12700 EMU.
$EMU = 12700.0

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Add-SubscriptLabel {
    param(
        [double]$OffX,
        [double]$OffY,
        [double]$ExtCx,
        [double]$ExtCy,
        [string]$ShapeName,
        [string]$MainText,
        [string]$SubText
    )

    $left   = $OffX / $EMU
    $top    = $OffY / $EMU
    $width  = $ExtCx / $EMU
    $height = $ExtCy / $EMU

    $tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
    $tb.Name = $ShapeName

    # <a:noFill/>
    $tb.Fill.Visible = $false

    # <a:bodyPr wrap="square"><a:spAutoFit/></a:bodyPr>
    $tb.TextFrame.WordWrap = $true
    $tb.TextFrame.AutoSize = 1

    $tf = $tb.TextFrame
    $tr = $tf.TextRange
    $tr.Text = $MainText + $SubText

    $mainLen = $MainText.Length
    $subLen = $SubText.Length

    $r1 = $tr.Characters(1, $mainLen)
    $r1.Font.Name = "Times"
    $r1.Font.NameFarEast = "Times"
    $r1.Font.NameComplexScript = "Times"
    $r1.Font.Size = 18
    $r1.Font.Italic = $true

    $r2 = $tr.Characters($mainLen + 1, $subLen)
    $r2.Font.Name = "Times"
    $r2.Font.NameFarEast = "Times"
    $r2.Font.NameComplexScript = "Times"
    $r2.Font.Size = 18
    $r2.Font.Italic = $true
    $r2.Font.Subscript = $true

    return $tb
}

# TextBox 60 (id 61): "l" + subscript "i"
Add-SubscriptLabel 2437766 438725 455579 369332 "TextBox 60" "l" "i" | Out-Null

# TextBox 61 (id 62): "l" + subscript "i+1"
Add-SubscriptLabel 2399845 3592028 642092 369332 "TextBox 61" "l" "i+1" | Out-Null

Write-Output "Added TextBox 60 and TextBox 61"
